$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(45, 2).Value = "79-367"
$ws.Cells.Item(46, 2).Value = "79-392"
$ws.Cells.Item(47, 2).Value = "79-393"
$ws.Cells.Item(48, 2).Value = "79-415"
$ws.Cells.Item(49, 2).Value = "79-419"
$ws.Cells.Item(50, 2).Value = "80-330"
$ws.Cells.Item(51, 2).Value = "82-180"
$ws.Cells.Item(52, 2).Value = "82-215"
$ws.Cells.Item(53, 2).Value = "82-234"
$ws.Cells.Item(54, 2).Value = "82-253"
$ws.Cells.Item(202, 2).Value = "60-142"
$ws.Cells.Item(203, 2).Value = "62-142"
$ws.Cells.Item(204, 2).Value = "70-350"
$ws.Cells.Item(205, 2).Value = "76-239"
$ws.Cells.Item(206, 2).Value = "76-260"
$ws.Cells.Item(207, 2).Value = "76-261"
$ws.Cells.Item(208, 2).Value = "76-265"
$ws.Cells.Item(209, 2).Value = "76-278"
$ws.Cells.Item(210, 2).Value = "79-317"
$ws.Cells.Item(211, 2).Value = "79-345"
$ws.Cells.Item(212, 2).Value = "79-465"
$ws.Cells.Item(213, 2).Value = "82-119"
$ws.Cells.Item(214, 2).Value = "82-130"
$ws.Cells.Item(215, 2).Value = "82-137"
$ws.Cells.Item(216, 2).Value = "82-237"
$ws.Cells.Item(217, 2).Value = "82-255"
$ws.Cells.Item(218, 2).Value = "82-278"
$ws.Cells.Item(219, 2).Value = "82-279"
$ws.Cells.Item(220, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---The Arts"
$ws.Cells.Item(220, 2).Value = "82-285"
$ws.Cells.Item(221, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---The Arts"
$ws.Cells.Item(221, 2).Value = "82-380"
$ws.Cells.Item(222, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---The Arts"
$ws.Cells.Item(222, 2).Value = "D7-001"
$ws.Cells.Item(223, 2).Value = "15-155"
$ws.Cells.Item(224, 2).Value = "21-102"
$ws.Cells.Item(225, 2).Value = "21-127"
$ws.Cells.Item(226, 2).Value = "73-155"
$ws.Cells.Item(227, 2).Value = "73-347"
$ws.Cells.Item(228, 2).Value = "79-175"
$ws.Cells.Item(229, 2).Value = "80-210"
$ws.Cells.Item(230, 2).Value = "80-150"
$ws.Cells.Item(231, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---Logic/Mathematical Reasoning"
$ws.Cells.Item(231, 2).Value = "80-211"
$ws.Cells.Item(232, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---Logic/Mathematical Reasoning"
$ws.Cells.Item(232, 2).Value = "80-312"
$ws.Cells.Item(233, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---Logic/Mathematical Reasoning"
$ws.Cells.Item(233, 2).Value = "D8-001"
$ws.Cells.Item(234, 2).Value = "05-360"
$ws.Cells.Item(235, 2).Value = "17-313"
$ws.Cells.Item(236, 2).Value = "19-101"
$ws.Cells.Item(237, 2).Value = "36-315"
$ws.Cells.Item(238, 2).Value = "48-095"
$ws.Cells.Item(239, 2).Value = "49-101"
$ws.Cells.Item(240, 2).Value = "70-100"
$ws.Cells.Item(241, 2).Value = "70-122"
$ws.Cells.Item(242, 2).Value = "70-246"
$ws.Cells.Item(243, 2).Value = "70-415"
$ws.Cells.Item(244, 2).Value = "76-270"
$ws.Cells.Item(245, 2).Value = "85-107"
$ws.Cells.Item(246, 2).Value = "88-150"
$ws.Cells.Item(247, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---Additional Disciplines (Business/Design/Engineering)"
$ws.Cells.Item(247, 2).Value = "88-223"
$ws.Cells.Item(248, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---Additional Disciplines (Business/Design/Engineering)"
$ws.Cells.Item(248, 2).Value = "88-235"
$ws.Cells.Item(249, 1).Value = "GenEd---GenEd---Disciplinary Perspectives---Additional Disciplines (Business/Design/Engineering)"
$ws.Cells.Item(249, 2).Value = "D9-001"
$ws.Cells.Item(250, 2).Value = "66-122"
$ws.Cells.Item(251, 2).Value = "66-132"
$ws.Cells.Item(252, 2).Value = "66-134"
$ws.Cells.Item(253, 2).Value = "66-138"
$ws.Cells.Item(254, 2).Value = "66-146"
$ws.Cells.Item(255, 2).Value = "66-147"
$ws.Cells.Item(256, 1).Value = "GenEd---GenEd---Special Seminars---Grand Challenge Seminar"
$ws.Cells.Item(256, 2).Value = "66-151"
$ws.Cells.Item(257, 1).Value = "GenEd---GenEd---Special Seminars---Grand Challenge Seminar"
$ws.Cells.Item(257, 2).Value = "66-181"
$ws.Cells.Item(258, 1).Value = "GenEd---GenEd---Special Seminars---Grand Challenge Seminar"
$ws.Cells.Item(258, 2).Value = "66-182"
$ws.Cells.Item(259, 2).Value = "36-301"
$ws.Cells.Item(260, 2).Value = "66-236"
$ws.Cells.Item(261, 2).Value = "76-219"
$ws.Cells.Item(262, 2).Value = "76-327"
$ws.Cells.Item(263, 2).Value = "76-442"
$ws.Cells.Item(264, 2).Value = "79-210"
$ws.Cells.Item(265, 2).Value = "79-215"
$ws.Cells.Item(266, 2).Value = "79-237"
$ws.Cells.Item(267, 2).Value = "79-270"
$ws.Cells.Item(268, 2).Value = "79-276"
$ws.Cells.Item(269, 2).Value = "79-313"
$ws.Cells.Item(270, 2).Value = "79-321"
$ws.Cells.Item(271, 2).Value = "79-349"
$ws.Cells.Item(272, 2).Value = "79-360"
$ws.Cells.Item(273, 2).Value = "79-378"
$ws.Cells.Item(274, 2).Value = "79-380"
$ws.Cells.Item(275, 2).Value = "79-432"
$ws.Cells.Item(276, 2).Value = "80-234"
$ws.Cells.Item(277, 2).Value = "80-245"
$ws.Cells.Item(278, 2).Value = "80-324"
$ws.Cells.Item(279, 2).Value = "80-334"
$ws.Cells.Item(280, 2).Value = "80-335"
$ws.Cells.Item(281, 2).Value = "80-348"
$ws.Cells.Item(282, 2).Value = "82-184"
$ws.Cells.Item(283, 2).Value = "82-224"
$ws.Cells.Item(284, 2).Value = "82-260"
$ws.Cells.Item(285, 2).Value = "82-289"
$ws.Cells.Item(286, 2).Value = "82-299"
$ws.Cells.Item(287, 2).Value = "82-399"
$ws.Cells.Item(288, 2).Value = "84-309"
$ws.Cells.Item(289, 2).Value = "84-322"
$ws.Cells.Item(290, 1).Value = "GenEd---GenEd---Special Seminars---Perspectives on Justice and Injustice"
$ws.Cells.Item(290, 2).Value = "84-352"
$ws.Cells.Item(290, 3).Value = 9
$ws.Cells.Item(291, 1).Value = "GenEd---GenEd---Special Seminars---Perspectives on Justice and Injustice"
$ws.Cells.Item(291, 2).Value = "84-367"
$ws.Cells.Item(291, 3).Value = 9
$ws.Cells.Item(292, 1).Value = "GenEd---GenEd---Special Seminars---Perspectives on Justice and Injustice"
$ws.Cells.Item(292, 2).Value = "85-350"
$ws.Cells.Item(292, 3).Value = 9
$ws.Cells.Item(293, 1).Value = "GenEd---GenEd---Special Seminars---Perspectives on Justice and Injustice"
$ws.Cells.Item(293, 2).Value = "88-221"
$ws.Cells.Item(293, 3).Value = 9
$ws.Cells.Item(294, 1).Value = "GenEd---GenEd---Special Seminars---Perspectives on Justice and Injustice"
$ws.Cells.Item(294, 2).Value = "88-285"
$ws.Cells.Item(294, 3).Value = 9
$ws.Cells.Item(295, 1).Value = "GenEd---GenEd---Special Seminars---Perspectives on Justice and Injustice"
$ws.Cells.Item(295, 2).Value = "99-262"
$ws.Cells.Item(295, 3).Value = 9
$ws.Cells.Item(296, 1).Value = "GenEd---GenEd---Special Seminars---Perspectives on Justice and Injustice"
$ws.Cells.Item(296, 2).Value = "99-265"
$ws.Cells.Item(296, 3).Value = 9
